# TRIMAZKON_task settings: renumber tasks and add the new task_2 / task_3 rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("task_settings")

$path = "C:/Users/jakub.hlavacek.local/Desktop/JHV/test_images/Keyence/_503_Witte/datumovka/A/Height_test/"

# Row 1 - TRIMAZKON_task_2 (was task_4's old slot, now holds task_2's data)
$ws.Cells.Item(1, 1).Value = "TRIMAZKON_task_2"
$ws.Cells.Item(1, 2).Value = $path
$ws.Cells.Item(1, 3).Value = "278"
$ws.Cells.Item(1, 4).Value = "998"
$ws.Cells.Item(1, 5).Value = "12:00"
$ws.Cells.Item(1, 6).ClearContents()

# Row 2 - TRIMAZKON_task_3 (new row)
$ws.Cells.Item(2, 1).Value = "TRIMAZKON_task_3"
$ws.Cells.Item(2, 2).Value = $path
$ws.Cells.Item(2, 3).Value = "200"
$ws.Cells.Item(2, 4).Value = "111"
$ws.Cells.Item(2, 5).Value = "12:00"

# Row 3 - TRIMAZKON_task_1 (unchanged data, shifted down one row)
$ws.Cells.Item(3, 1).Value = "TRIMAZKON_task_1"
$ws.Cells.Item(3, 2).Value = $path
$ws.Cells.Item(3, 3).Value = "30"
$ws.Cells.Item(3, 4).Value = "500"
$ws.Cells.Item(3, 5).Value = "23:59"

# Row 4 - TRIMAZKON_task_4 (new row with task_4's new data)
$ws.Cells.Item(4, 1).Value = "TRIMAZKON_task_4"
$ws.Cells.Item(4, 2).Value = $path
$ws.Cells.Item(4, 3).Value = "278"
$ws.Cells.Item(4, 4).Value = "998"
$ws.Cells.Item(4, 5).Value = "1:00"
